$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 7-11, columns B:G with new forecast-error values
$data = @{
    7  = @(0.06372808675650223, 0.9343998086003413, 2.775446151470913, 1.665967031927977, 1.688362315275229, 36)
    8  = @(0.0399851282911067, 1.050477138034481, 3.029802060013869, 1.740632660848885, 1.765578673122826, 35)
    9  = @(-0.04534382660436193, 1.253794741722951, 4.508050214875826, 2.123216949554573, 2.177877806311915, 20)
    10 = @(-0.4013436200723844, 0.9917379900131201, 2.080463190564922, 1.442381083682437, 1.441990070455223, 13)
    11 = @(0.2679351619938605, 0.5450213326915778, 0.3698384655848825, 0.6081434580630483, 0.6103781763712214, 5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        # Column B is index 2, C is 3, D is 4, E is 5, F is 6, G is 7
        $col = $i + 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
